$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 onto the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Innings-pitched detail values (I0 = innings at outing start, IF = innings at outing finish)
$data = @{
    2 = @(6, 6)
    3 = @(8, 8)
    4 = @(8, 8)
    5 = @(8, 8)
    6 = @(7, 8)
    7 = @(8, 8)
    8 = @(8, 8)
    9 = @(9, 9)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(7, 7)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(7, 8)
    18 = @(7, 7)
    19 = @(8, 8)
    20 = @(10, 10)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(7, 7)
    24 = @(8, 9)
    25 = @(6, 6)
    26 = @(9, 9)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(6, 7)
    30 = @(6, 7)
    31 = @(8, 8)
    32 = @(8, 9)
    33 = @(7, 7)
    34 = @(10, 10)
    35 = @(6, 7)
    36 = @(7, 8)
    37 = @(8, 8)
    38 = @(6, 7)
    39 = @(8, 8)
    40 = @(8, 8)
    41 = @(7, 7)
    42 = @(7, 8)
    43 = @(9, 9)
    44 = @(7, 8)
    45 = @(7, 7)
    46 = @(7, 7)
    47 = @(6, 6)
    48 = @(6, 6)
    49 = @(6, 7)
    50 = @(7, 7)
    51 = @(6, 7)
    52 = @(7, 7)
    53 = @(7, 7)
    54 = @(4, 5)
    55 = @(7, 7)
    56 = @(6, 7)
    57 = @(6, 6)
    58 = @(6, 6)
    59 = @(1, 1)
    60 = @(6, 6)
    61 = @(4, 4)
    62 = @(4, 4)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Host "I0/IF columns populated"
